$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting existing rows 3-12 down to 4-13
$ws.Rows("3:3").Insert()

# Set the new row 3 values. Most columns (A,B,C,E,F,G,H,I,J,K,L,Q,T) are same as row 4 (the shifted-down old row 3)
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 45107
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107001
$ws.Range("J3").Value = "Caqui"
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

# Apply the date cell number format to D3 to match D4 (date number format)
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
